$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so the exact string representation is preserved
$textCells = @("D5", "D6", "D10", "D15", "D19", "D21", "D22", "D24", "D27", "D31", "D32", "D33", "D36", "D37", "D38", "D39", "D42", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "72.258.35"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.638.23"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "587.74"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "175.10"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").Value = "2.637.56"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "3.121.20"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "72.168.30"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.641.05"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").Value = "12.10"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "373.90"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -1.39%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").Value = "71.42"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").Value = "2.774.86"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "7.94"
$ws.Range("E31").Value = "  -1.12%  "
$ws.Range("D32").Value = "492.00"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "161.30"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +8.48%  "
$ws.Range("D38").Value = "19.15"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "18.90"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "1.73"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "0.326"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").Value = "39.03"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "150.28"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "3.63"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").Value = "0.605"
$ws.Range("E51").Value = "  +1.10%  "
